# "Corrijo costo del logger"
# The logger's battery ("Pilas") quantity on the Costos sheet (cell H13)
# was wrong; fix it from 1 to 2. All dependent formulas (I13, I18, and the
# downstream cash-flow figures on "Flujo de caja") recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsCostos = $wb.Worksheets.Item("Costos")
$wsCostos.Range("H13").Value = 2

# Also fix a leftover copy/paste label on "Ingresos": the second block's
# header ("Bases") still read "Loggers" (copied from the first block).
$wsIngresos = $wb.Worksheets.Item("Ingresos")
$wsIngresos.Range("D6").Value = "Bases"

# Leave the workbook focused back on "Costos" (where the fix was made),
# matching the editor's final on-screen state.
$wsCostos.Activate()
$wsCostos.Range("H22").Select()
